# Applies the commit "idk measurements and stuff":
#  - adds two new worksheets: 2_stick_measurements, 3_particle_measurements
#  - populates them with stick/particle diameter measurement data
#  - changes the active sheet/tab selection + the 0_raw_laser selection

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 0_raw_laser
$ws2 = $wb.Worksheets.Item(2)   # 1_filtered_data

# --- add the two new worksheets at the end, in order ---------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $last)
$ws3.Name = "2_stick_measurements"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $last)
$ws4.Name = "3_particle_measurements"

# ===========================================================================
# Sheet 3: 2_stick_measurements
# ===========================================================================
# Fill data cells first (this is the order the shared-string table records
# these values in: the picture names for the stick-diameter shots, then the
# four new headers typed from column F back to column C).
$stickNames = @("DSC03209.JPG","DSC03210.JPG","DSC03211.JPG","DSC03212.JPG","DSC03213.JPG","DSC03214.JPG")
$hDiff      = @(2105,2112,2112,2099,2093,2106)
$vDiff      = @(131,111,48,84,36,71)

for ($i = 0; $i -lt $stickNames.Length; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r, 2).Value = $stickNames[$i]
}

$ws3.Range("F1").Value = "diameter_length_error-pixels"
$ws3.Range("E1").Value = "diameter_length-pixels"
$ws3.Range("C1").Value = "diameter_measure_horizontal_diff"
$ws3.Range("D1").Value = "diameter_measure_vertical_diff"
$ws3.Range("B1").Value = "picture_name"

for ($i = 0; $i -lt $stickNames.Length; $i++) {
    $r = $i + 2
    $ws3.Cells.Item($r, 1).Value = $i
    $ws3.Cells.Item($r, 3).Value = $hDiff[$i]
    $ws3.Cells.Item($r, 4).Value = $vDiff[$i]
    $ws3.Cells.Item($r, 5).Formula = "=SQRT(C$r^2+D$r^2)"
    $ws3.Cells.Item($r, 6).Value = 5
}

# --- formatting: reuse the same header / index styles as the other sheets -
$ws1.Range("B1").Copy()
$ws3.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws3.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ===========================================================================
# Sheet 4: 3_particle_measurements
# ===========================================================================
$ws4.Range("B2").Value = "DSC03162.JPG"
$ws4.Range("B3").Value = "DSC03163.JPG"
$ws4.Range("C1").Value = "diameter-pixels"
$ws4.Range("B1").Value = "picture_name"
$ws4.Range("B4").Value = "DSC03175.JPG"

$ws4.Range("A2").Value = 0
$ws4.Range("A3").Value = 1
$ws4.Range("A4").Value = 2
$ws4.Range("C2").Value = 48
$ws4.Range("C3").Value = 49
$ws4.Range("C4").Value = 52

$ws1.Range("B1").Copy()
$ws4.Range("B1:C1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws4.Range("A2:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ===========================================================================
# Selections / active sheet
# ===========================================================================
$ws3.Range("A1:C4").Select()
$ws1.Range("A1:D6").Select()
$ws4.Range("C1").Select()
$ws4.Activate()
